$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / non-numeric-looking values: direct assignment is safe
$ws.Range("D2").Value = "34.784.96"
$ws.Range("E2").Value = "  -1.75%  "
$ws.Range("D3").Value = "1.812.27"
$ws.Range("E3").Value = "  -1.34%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("E5").Value = "  +0.27%  "
$ws.Range("E6").Value = "  -0.59%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").Value = "  -8.28%  "
$ws.Range("E9").Value = "  +5.44%  "
$ws.Range("E10").Value = "  -2.42%  "
$ws.Range("E11").Value = "  -1.58%  "
$ws.Range("D12").Value = "2.072.58"
$ws.Range("E12").Value = "  -1.42%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("E13").Value = "  +0.76%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("E14").Value = "  -0.43%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.795.23"
$ws.Range("E15").Value = "  -2.28%  "
$ws.Range("E16").Value = "  -1.38%  "
$ws.Range("D17").Value = "34.744.39"
$ws.Range("E17").Value = "  -1.74%  "
$ws.Range("E18").Value = "  -0.65%  "
$ws.Range("D19").Value = "0.0₃0785"
$ws.Range("E19").Value = "  -1.52%  "
$ws.Range("E20").Value = "  -1.45%  "
$ws.Range("E21").Value = "  -0.74%  "
$ws.Range("E22").Value = "  +0.71%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("E24").Value = "  +2.20%  "
$ws.Range("E25").Value = "  +1.55%  "
$ws.Range("E26").Value = "  -1.51%  "
$ws.Range("E27").Value = "  -2.52%  "
$ws.Range("E28").Value = "  -0.72%  "
$ws.Range("E29").Value = "  +0.19%  "
$ws.Range("E30").Value = "  +0.21%  "
$ws.Range("E31").Value = "  +3.62%  "
$ws.Range("E32").Value = "  -0.50%  "
$ws.Range("E33").Value = "  -2.28%  "
$ws.Range("E34").Value = "  +16.27%  "
$ws.Range("E35").Value = "  -3.18%  "
$ws.Range("E36").Value = "  +2.45%  "
$ws.Range("E37").Value = "  -4.23%  "
$ws.Range("D39").Value = "1.332.44"
$ws.Range("E39").Value = "  -0.87%  "
$ws.Range("E40").Value = "  -0.89%  "
$ws.Range("E41").Value = "  +0.11%  "
$ws.Range("E42").Value = "  -3.94%  "
$ws.Range("E43").Value = "  -6.87%  "
$ws.Range("E44").Value = "  -8.39%  "
$ws.Range("E45").Value = "  -4.52%  "
$ws.Range("E46").Value = "  +0.99%  "
$ws.Range("E47").Value = "  -1.26%  "
$ws.Range("D48").Value = "1.999.65"
$ws.Range("E48").Value = "  -0.23%  "
$ws.Range("E49").Value = "  +0.10%  "
$ws.Range("E50").Value = "  +7.29%  "
$ws.Range("E51").Value = "  -4.38%  "

# Numeric-looking price strings: force text storage so trailing/insignificant
# zeros and multi-dot "thousand.thousand.cents" formatting survive, matching
# the source feed's plain-text price column, then clear the temporary text
# number-format back off so no stray style is left behind on the cell.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.11"
$ws.Range("D5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.605"
$ws.Range("D6").ClearFormats()
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "39.81"
$ws.Range("D8").ClearFormats()
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0682"
$ws.Range("D10").ClearFormats()
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0993"
$ws.Range("D11").ClearFormats()
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.29"
$ws.Range("D13").ClearFormats()
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.667"
$ws.Range("D14").ClearFormats()
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.60"
$ws.Range("D16").ClearFormats()
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.49"
$ws.Range("D18").ClearFormats()
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "241.11"
$ws.Range("D20").ClearFormats()
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.93"
$ws.Range("D21").ClearFormats()
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.70"
$ws.Range("D22").ClearFormats()
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "171.81"
$ws.Range("D25").ClearFormats()
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.76"
$ws.Range("D26").ClearFormats()
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.25"
$ws.Range("D27").ClearFormats()
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.120"
$ws.Range("D28").ClearFormats()
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.06"
$ws.Range("D31").ClearFormats()
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0547"
$ws.Range("D32").ClearFormats()
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.96"
$ws.Range("D33").ClearFormats()
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.26"
$ws.Range("D34").ClearFormats()
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.700"
$ws.Range("D36").ClearFormats()
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "92.02"
$ws.Range("D37").ClearFormats()
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.967"
$ws.Range("D42").ClearFormats()
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "14.32"
$ws.Range("D43").ClearFormats()
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.22"
$ws.Range("D44").ClearFormats()
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.27"
$ws.Range("D46").ClearFormats()
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "97.97"
$ws.Range("D51").ClearFormats()
